$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-43) holds a date serial that was bumped from 45780 to 45781
# (2025-05-03 -> 2025-05-04). Update the value while leaving the existing
# date-formatted style (s="1") untouched.
$ws.Range("C2:C43").Value = 45781
